$d = $word.ActiveDocument

function Replace-All($old, $new) {
  $null = $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}

function Replace-InParagraph($paraIndex, $old, $new) {
  $p = $d.Paragraphs($paraIndex)
  $r = $p.Range
  $null = $r.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}

# 'Thai' appears twice with two different French translations:
# nav-list link -> 'Thailandais', section heading -> 'Thailande'.
# Disambiguate using the distinct paragraphs that contain them.
$thaiDone = $false
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
  $p = $d.Paragraphs($i)
  $txt = $p.Range.Text
  if (-not $thaiDone -and $txt -match '^English / Portuguese / French / Thai / Vietnamese / Spanish') {
    Replace-InParagraph $i 'Thai' 'Thaïlandais'
    $thaiDone = $true
  }
}
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
  $p = $d.Paragraphs($i)
  $txt = $p.Range.Text
  if ($txt -eq "Thai`r") {
    Replace-InParagraph $i 'Thai' 'Thaïlande'
  }
}

# 'Subject line' -> 'Ligne de sujet' only in the English-section header;
# all other remaining occurrences -> "Ligne d'objet".
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
  $p = $d.Paragraphs($i)
  $txt = $p.Range.Text
  if ($txt -match '^Subject line: Meet our team') {
    Replace-InParagraph $i 'Subject line' 'Ligne de sujet'
  }
}
Replace-All 'Subject line' 'Ligne d''objet'

# ' or ' -> ' ou ' only the one immediately after the English 'live chat' link
# (a second, untouched ' or ' exists later in the same paragraph, and another
# one in the French section that must stay as-is).
$rngLC = $d.Content
$foundLC = $rngLC.Find.Execute('live chat', $true, $true, $false, $false, $false, $true, 1, $false, '', 0)
if ($foundLC) {
  $narrow = $d.Range($rngLC.End, $rngLC.End + 10)
  $null = $narrow.Find.Execute(' or ', $false, $true, $false, $false, $false, $true, 1, $false, ' ou ', 1)
}

# ' o ' -> ' ou ' only the one immediately after the Spanish 'live chat' link.
$rngLC2 = $d.Content
$foundLC2 = $rngLC2.Find.Execute('live chat', $true, $true, $false, $false, $false, $true, 1, $false, '', 0)
while ($foundLC2) {
  $checkRange = $d.Range($rngLC2.End, [Math]::Min($rngLC2.End + 3, $d.Content.End))
  if ($checkRange.Text -eq ' o ') {
    $narrow2 = $d.Range($rngLC2.End, $rngLC2.End + 5)
    $null = $narrow2.Find.Execute(' o ', $false, $true, $false, $false, $false, $true, 1, $false, ' ou ', 1)
    $foundLC2 = $false
  } else {
    $foundLC2 = $rngLC2.Find.Execute('live chat', $true, $true, $false, $false, $false, $true, 1, $false, '', 0)
  }
}

# Remaining, unambiguous replacements (unique text or consistent across all
# duplicate occurrences).
Replace-All 'English' 'Anglais'
Replace-All 'Portuguese' 'Portugais'
Replace-All 'French' 'Français'
Replace-All 'Vietnamese' 'Vietnamien'
Replace-All 'Spanish' 'Espagnol'
Replace-All 'Brief' 'Bref'
Replace-All 'An email to partners in the the target country to invite them for a one-day seminar. It will be sent via customer.io' 'Un courriel adressé aux partenaires du pays cible pour les inviter à un séminaire d''une journée. Il sera envoyé via customer.io'
Replace-All 'Target audience' 'Public cible'
Replace-All 'Partners in the target country' 'Partenaires dans le pays cible'
Replace-All ': Meet our team in [CITY] | [DATE] ' ': Rencontrez notre équipe à [CITY] | [DATE] '
Replace-All 'You’re invited to our Deriv Partner Seminar' 'Vous êtes invités à notre séminaire Deriv Partner'
Replace-All 'Dear [PARTNER NAME], ' 'Cher [NOM DU PARTENAIRE], '
Replace-All 'We’re excited to let you know that the Deriv Affiliate team will be in [CITY] in [MONTH] to meet with you, our valued partners!' 'Nous sommes heureux de vous annoncer que l''équipe de Deriv Affiliate sera présente à [CITY] au mois de [MONTH] pour vous rencontrer, vous, nos précieux partenaires !'
Replace-All 'Your country manager will inform you about the exact location by [DATE]' 'Votre responsable de pays vous informera de l''emplacement exact à l''adresse suivante : [DATE].'
Replace-All 'In this one-day seminar, we’ll be providing technical and marketing support, offering the opportunity to network with other partners over a delicious lunch as well as listening to your feedback about our partnership programmes. This is your chance to get your voice heard, which will help us plan future efforts to support you better. ' 'Au cours de ce séminaire d''une journée, nous vous fournirons une assistance technique et marketing. Vous pourrez faire du réseautage avec d''autres partenaires autour d''un délicieux déjeuner, tout en nous faisant part de vos impressions sur nos programmes de partenariat. Il s''agit d''une occasion de faire entendre votre voix, ce qui nous aidera à prendre des initiatives à l''avenir pour mieux vous soutenir. '
Replace-All 'Please RSVP by submitting the registration form by ' 'Veuillez RSVP en envoyant le formulaire d''inscription avant '
Replace-All '. Please note that attendance is confirmed on a first come, first served basis. We look forward to seeing you there!' '. Veuillez noter que la participation est confirmée sur la base du premier arrivé, premier servi. Nous sommes impatients de vous y rencontrer !'
Replace-All 'Send my details' 'Envoyer mes coordonnées'
Replace-All 'If you have any questions, please contact us via ' 'Si vous avez des questions, veuillez nous contacter via '
Replace-All '. / If you have any questions, please contact your country manager, [NAME], at [EMAIL ADDRESS] or [WHATSAPP NO] (WhatsApp). ' '. / Si vous avez des questions, veuillez contacter votre responsable national, [NAME], à [ADRESSE ÉLECTRONIQUE] ou [N° WHATSAPP] (WhatsApp). '
Replace-All 'Back to' 'Retour à'
Replace-All 'Você está convidado(a) para o nosso ' 'Vous êtes convié(e) à nos réunions. '
Replace-All 'Prazado(a) [PARTNER NAME], ' 'Prazado(a) [NOM DU PARTENAIRE], '
Replace-All 'É com grande satisfação que comunicamos que a equipe de Afiliados da Deriv estará em [CITY] em [MONTH] para se reunir com vocês, nossos queridos parceiros!' 'C''est avec une grande satisfaction que nous vous annonçons qu''une équipe d''Afiliados da Deriv est sur le site [CITY] à [MONTH] pour se réunir avec vous, nos chers collègues !'
Replace-All 'Neste seminário de 1 dia, iremos oferecer suporte técnico e de marketing, oportunidades de interação com outros parceiros durante um excelente almoço, e também ouviremos os seus comentários sobre os nossos programas de parceria. Esta é sua grande oportunidade de fazer com que sua voz seja ouvida, para assim planejarmos ações futuras capazes de proporcionar a você um suporte ainda melhor. ' 'Lors de ce séminaire d''un jour, nous vous offrirons un soutien technique et marketing, des possibilités d''échanges avec d''autres partenaires au cours d''un excellent déjeuner, et nous vous présenterons les commentaires de chacun sur nos programmes de partenariat. Il s''agit d''une grande opportunité pour vous de faire en sorte que votre voix s''exprime bien, afin que nous puissions planifier des actions futures capables de vous fournir un soutien toujours meilleur. '
Replace-All 'Por favor, confirme sua presença enviando o formulário de cadastro até o dia ' 'Veuillez confirmer votre présence en envoyant le formulaire de candidature au jour '
Replace-All '. Observe que a participação será confirmada por ordem de chegada. Contamos com a sua presença! ' '. Observez que la participation est confirmée par l''ordre de passage. Nous vous aidons dans votre présentation ! '
Replace-All 'Enviar meus dados' 'Envoyer mes messages'
Replace-All 'Em caso de dúvida, entre em contato conosco através do ' 'En cas de doute, entrez en contact avec nous par le biais de '
Replace-All ' ou pelo ' ' ou par '
Replace-All ' em nosso site. / Se você tiver alguma dúvida, fale com o gerente de seu país [NAME] em [EMAIL ADDRESS] ou [WHATSAPP NO] (Whatsapp).' ' sur notre site. / Si vous avez des doutes, contactez le responsable de votre pays [NAME] à l''adresse [EMAIL ADDRESS] ou [WHATSAPP NO] (Whatsapp).'
Replace-All 'Veuillez confirmer votre présence en soumettant le formulaire d''inscription avant le [DATE]. Veuillez noter que la participation est confirmée selon le principe du premier arrivé, premier servi. Nous avons hâte de vous voir là-bas!' 'Veuillez confirmer votre présence en soumettant le formulaire d''inscription avant le [DATE]. Veuillez noter que la participation est confirmée selon le principe du premier arrivé, premier servi. Nous avons hâte de vous voir là-bas !'
Replace-All 'Deriv Partner Seminar ของเรา' 'Séminaire Deriv Partner ของเรา'
Replace-All 'เรียน [PARTNER NAME] ' 'เรียน [NOM DU PARTENAIRE] '
Replace-All 'เรารู้สึกตื่นเต้นยินดีที่จะแจ้งให้คุณทราบว่า ทีมงานพันธมิตรของ Deriv จะไปเยือนที่ [CITY] ในเดือน [MONTH] เพื่อพบกับคุณซึ่งเป็นพันธมิตรที่มีคุณค่าของเรา!' 'เรารู้สึกตื่นเต้นยินดีที่จะแจ้งให้คุณทราบว่า ทีมงานพันธมิตรของ Deriv จะไปเยือนที่ [CITY] ในเดือน [MONTH] เพื่อพบกับคุณซึ่งเป็นพันธมิตรที่มีคุณค่าของเรา !'
Replace-All ' และเราหวังว่าจะได้พบคุณที่นั่น!' ' และเราหวังว่าจะได้พบคุณที่นั่น !'
Replace-All 'Xin chào [PARTNER NAME], ' 'Xin chào [NOM DU PARTENAIRE], '
Replace-All 'Chúng tôi rất vui mừng thông báo đến bạn rằng đội ngũ Tiếp thị liên kết của Deriv sẽ có mặt tại [CITY] vào [MONTH] để gặp bạn, một đối tác quan trọng của chúng tôi!' 'Chúng tôi rất vui mừng thông báo đến bạn rằng đội ngũ Tiếp thị liên kết của Deriv sẽ có mặt tại [CITY] vào [MONTH] để gặp bạn, một đối tác quan trọng của chúng tôi !'
Replace-All 'Giám đốc phụ trách quốc gia sẽ thông báo địa điểm chính xác đến bạn trước ngày [DATE]' 'Giám đốc phụ trách quốc gia sẽ thông báo địa điểm chính xác đạn bạn trước ngày [DATE]'
Replace-All '. Xin lưu ý rằng chúng tôi sẽ ưu tiên xác nhận tham gia đối với những đơn đăng ký được gửi về trước. Chúng tôi mong được gặp bạn tại hội thảo!' '. Xin lưu ý rằng chúng tôi sẽ ưu tiên xác nhận tham gia đối với những đơn đăng ký được gửi về trước. Chúng tôi mong được gặp bạn tại hội thảo !'
Replace-All 'Te invitamos a nuestro seminario Deriv para socios' 'Vous êtes invités à notre séminaire Deriv pour les entreprises'
Replace-All 'Nos complace informarte que el equipo de Afiliados de Deriv estará en [CITY] en [MONTH] para reunirse con ustedes, ¡nuestros valiosos socios!' 'Nous avons le plaisir de vous informer que l''équipe d''Afiliados de Deriv se trouve à [CITY] en [MONTH] pour vous rencontrer, vous qui êtes de précieux partenaires !'
Replace-All 'Tu gerente de cuenta te informará de la ubicación exacta el [DATE]' 'Votre gérant de compte vous informera de l''emplacement exact le [DATE]'
Replace-All 'En este seminario de un día te brindaremos apoyo técnico y de marketing, ofreceremos la oportunidad de establecer contactos con otros socios durante un delicioso almuerzo y escucharemos tus comentarios sobre nuestros programas de asociación. Esta es tu oportunidad de hacer oír tu voz, que nos ayudará a planificar futuros esfuerzos para apoyarte mejor. ' 'Au cours de ce séminaire d''un jour, vous bénéficierez d''une assistance technique et marketing, vous aurez l''occasion d''établir des contacts avec d''autres membres de l''association au cours d''un délicieux repas et vous entendrez vos commentaires sur nos programmes d''association. C''est l''occasion de faire entendre votre voix, qui nous aidera à planifier nos futurs efforts pour mieux vous aider. '
Replace-All 'Por favor, confirma tu asistencia enviando el formulario de inscripción antes del ' 'Veuillez confirmer votre participation en envoyant le formulaire d''inscription avant le '
Replace-All '. Ten en cuenta que la asistencia se confirmará por orden de llegada. Esperamos verte allí.' '. Tenez compte du fait que l''assistance est confirmée par l''ordre d''arrivée. Esperamos verte allí.'
Replace-All 'Enviar mis datos' 'Envoi de vos données'
Replace-All 'Si tienes alguna pregunta, contáctanos a través del ' 'Si vous avez une question, contactez-nous par le biais de '
Replace-All ' en nuestro sitio web. / Si tienes alguna pregunta, contacta con el gerente de cuenta de tu país [NAME] en [EMAIL ADDRESS] o [WHATSAPP NO] (WhatsApp). ' ' sur notre site web. / Si vous avez des questions, contactez le responsable de la clientèle de votre pays [NAME] à [ADRESSE ÉLECTRONIQUE] ou [N° WHATSAPP] (WhatsApp). '
